# Updates cryptos list values (Price in column D, Volume(1h) in column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.553.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5246"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3967"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09681"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.536"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("E13").Value = "  +3.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.910.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.27%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.60%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.348"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.638.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.296"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.693"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.130.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1086"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.753"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.913"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06753"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02437"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.263"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2226"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.092"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6468"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.191"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +2.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6097"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.285"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("E49").Value = "  +4.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("E51").Value = "  +2.49%  "

